$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Content.Find.Execute("2024-06-28 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-29 Saturday", 2) | Out-Null

# Update each arithmetic problem cell in the table, by position, to avoid
# ambiguity from duplicate problem text (e.g. "43-41=" appears twice with
# different replacements depending on position).
$newValues = @(
    "31-20=",
    "66-36=",
    "91-42=",
    "23+61=",
    "83-34=",
    "57-37=",
    "21-1=",
    "87-11=",
    "80+16=",
    "65+21=",
    "77-19=",
    "96-81=",
    "11+88=",
    "6+22=",
    "80-15=",
    "49+23=",
    "25+18=",
    "74-48=",
    "68-19=",
    "66-65=",
    "82-31=",
    "51-24=",
    "34-19=",
    "49+44=",
    "39-38=",
    "92-82=",
    "42+9=",
    "26+70=",
    "58-21=",
    "85-35=",
    "44-16=",
    "60+13=",
    "58-30=",
    "50+30=",
    "64-40=",
    "36+53=",
    "19-8=",
    "20+59=",
    "48-34=",
    "67-0=",
    "8+64=",
    "59-14=",
    "58-44=",
    "17+12=",
    "77-56=",
    "33+31=",
    "28-2=",
    "81-43=",
    "29+21=",
    "72-45=",
    "47+16=",
    "23-20=",
    "13-8=",
    "31-26=",
    "87-28=",
    "69-11=",
    "41-17=",
    "66+28=",
    "49+30=",
    "36+5=",
    "26+64=",
    "85+4=",
    "99-86=",
    "70-64=",
    "69+18=",
    "96-95=",
    "1+8=",
    "91-60=",
    "6-4=",
    "89-56=",
    "46+36=",
    "44+44=",
    "60+5=",
    "15+12=",
    "76-22=",
    "76-70=",
    "54-10=",
    "11+44=",
    "65-42=",
    "5+8=",
    "24+43=",
    "91-24=",
    "31+0=",
    "54-20=",
    "48-42=",
    "60-11=",
    "8+77=",
    "6+15=",
    "36-2=",
    "79-35=",
    "29+65=",
    "53-16=",
    "2+22=",
    "25+26=",
    "49+9=",
    "89-87=",
    "67+21=",
    "28+7=",
    "66+8=",
    "11+34="
)

$t = $d.Tables.Item(1)
$cols = 5
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        $rng.MoveEnd(1, -1) | Out-Null
        $rng.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated" $idx "cells"